# LeakSensor_BOM.xlsx update:
# The CON4 header connector line item now reflects a quantity of 1 (instead of 3),
# so its extended Total drops to 0.65 (1 x $0.65) instead of 1.95 (3 x $0.65).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.65

# Leave the cursor parked on G4, matching the saved selection state.
$ws.Range("G4").Select()
